# Update forest data - 2026-02-16 12:33
#
# 1) The two rows currently on the "New" sheet (rows 2-3) are promoted /
#    moved onto the end of the "Previously added" sheet (as rows 500-501),
#    keeping their values, hyperlinks and styling intact.
# 2) Four brand-new listings take their place as rows 2-5 on the "New"
#    sheet, each with its own hyperlink.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Previously added")
$ws2 = $wb.Worksheets.Item("New")

# ---------------------------------------------------------------------
# Step 0: stash a pristine copy of the "New" sheet's row styling
# (s=3/4/4/4/4/2) off in unused scratch rows far below the data, so we
# have an untouched formatting template to paste from later even after
# rows 2-5 themselves have been overwritten with new values.
# ---------------------------------------------------------------------
$ws2.Range("A3:F3").Copy($ws2.Range("A1000"))
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# Step 1: move the existing "New" rows (2 and 3) to the bottom of
# "Previously added" as rows 500 and 501.
# ---------------------------------------------------------------------

# Capture the hyperlink targets before we touch anything. Throughout this
# workbook the hyperlink target is always identical to the cell's own
# displayed URL text, so read it straight from the cell.
$oldUrl2 = $ws2.Range("A2").Value()
$oldUrl3 = $ws2.Range("A3").Value()

# Pre-register the hyperlinks on the destination cells first so that the
# subsequent full-value/format copy (which lays down the correct
# hyperlink-text style) has the final say on cell styling.
$ws1.Hyperlinks.Add($ws1.Range("A500"), $oldUrl2)
$ws1.Hyperlinks.Add($ws1.Range("A501"), $oldUrl3)

# Copy values + formats from the "New" sheet rows into place.
$ws2.Range("A2:F3").Copy($ws1.Range("A500"))

# ---------------------------------------------------------------------
# Step 2: replace rows 2-5 on "New" with four fresh rows, all styled
# like the scratch template row (s=3/4/4/4/4/2).
# ---------------------------------------------------------------------

# Drop every existing hyperlink on "New" - they'll all be re-created below.
$ws2.Hyperlinks.Delete()

$ws2.Range("A1000:F1000").Copy()
$ws2.Range("A2:F2").PasteSpecial(-4122)
$ws2.Range("A3:F3").PasteSpecial(-4122)
$ws2.Range("A4:F4").PasteSpecial(-4122)
$ws2.Range("A5:F5").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$rows = @(
    @{ Row = 2; A = "https://www.ss.com/msg/lv/real-estate/wood/daugavpils-and-reg/sventes-pag/ombfj.html";  B = "5 000 €";  C = "Daugavpils un raj."; D = "2 ha.";    E = "44880040024"; F = 46066.75833333333 },
    @{ Row = 3; A = "https://www.ss.com/msg/lv/real-estate/wood/madona-and-reg/sarkanu-pag/iegck.html";      B = "45 000 €"; C = "Madona un raj.";     D = "5 ha.";    E = "70900080055"; F = 46068.87291666667 },
    @{ Row = 4; A = "https://www.ss.com/msg/lv/real-estate/wood/preili-and-reg/rozkalnu-pag/hodnc.html";     B = "25 000 €"; C = "Preiļi un raj.";     D = "6.40 ha."; E = "76640040150"; F = 46069.334027777775 },
    @{ Row = 5; A = "https://www.ss.com/msg/lv/real-estate/wood/preili-and-reg/preilu-pag/plmkm.html";       B = "14 000 €"; C = "Preiļi un raj.";     D = "1.52 ha."; E = "76780040110"; F = 46068.884722222225 }
)

foreach ($r in $rows) {
    $rowNum = $r.Row
    $ws2.Range("A$rowNum").Value = $r.A
    $ws2.Range("B$rowNum").Value = $r.B
    $ws2.Range("C$rowNum").Value = $r.C
    $ws2.Range("D$rowNum").Value = $r.D
    # The E column values are digit-only strings; force them to stay text
    # (rather than being coerced to numbers) with a leading apostrophe.
    # This can reset the cell's style, so immediately restamp just that
    # one cell's formatting from the untouched scratch template.
    $ws2.Range("E$rowNum").Value = "'" + $r.E
    $ws2.Range("E1000").Copy()
    $ws2.Range("E$rowNum").PasteSpecial(-4122)
    $excel.CutCopyMode = $false
    $ws2.Range("F$rowNum").Value = $r.F
}

foreach ($r in $rows) {
    $rowNum = $r.Row
    $ws2.Hyperlinks.Add($ws2.Range("A$rowNum"), $r.A)
    # Adding the hyperlink re-styles column A on its own, so restamp it
    # from the untouched scratch template one more time.
    $ws2.Range("A1000").Copy()
    $ws2.Range("A$rowNum").PasteSpecial(-4122)
    $excel.CutCopyMode = $false
}

# Clean up the scratch template row.
$ws2.Range("A1000:F1000").Clear()

Write-Output "Forest data updated."
